$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header value updates
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2: move value from D2 to C2 (clear D2, set C2)
$ws.Range("D2").ClearContents()
$ws.Range("C2").Value = 42.682147309915706

# Row 3: clear B3 and C3
$ws.Range("B3").ClearContents()
$ws.Range("C3").ClearContents()

# Update selection to match new sqref B1:E3
$ws.Range("B1:E3").Select()
